# Fix sorting and generate viable xlsx and charts
#
# The two fastest runs for RandomInt33 / InsertBinarySortTimes.csv were
# re-measured; update the Avg_Time_ms values for the 5,000-row and
# 10,000-row samples (rows 2 and 3 of the Data sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2: Rows=5000  Avg_Time_ms 19.153686 -> 18.734492
$ws.Range("D2").Value = 18.734492

# Row 3: Rows=10000 Avg_Time_ms 73.476158 -> 76.22014
$ws.Range("D3").Value = 76.22014

# The worksheet cells are the chart's data source (Data!$D$2:$D$8), so the
# scatter chart picks up these corrected values automatically.
